$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column F ("Power") values for rows 2-20 were recorded in the wrong unit
# (kW instead of W) - multiply each by 1000 to correct them.
for ($row = 2; $row -le 20; $row++) {
    $cell = $ws.Cells.Item($row, 6)
    $cell.Value2 = $cell.Value2 * 1000
}
